$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: park_walk -- picture/description changed, link added
$ws.Range("J5").Value = "park_walk.jpg"
$ws.Range("K5").Value = "Walk and play the in park"
$ws.Range("L5").Value = "https://unsplash.com/photos/yIJIO2dhWWY"

# Row 6: ice_skating -- picture changed, link added
$ws.Range("J6").Value = "ice_skating.jpg"
$ws.Range("L6").Value = "https://unsplash.com/photos/lnCSMikKqfw"

# Row 7: beach_picnic -- picture changed, link added
$ws.Range("J7").Value = "beach_picnic.jpg"
$ws.Range("L7").Value = "https://unsplash.com/photos/hXY98KmQWkI"

# Row 8: beach_swimming -- picture changed, link added
$ws.Range("J8").Value = "beach_swimming.jpg"
$ws.Range("L8").Value = "https://unsplash.com/photos/R_BLOGXpsOg"

# Row 9: snow_angels -- picture changed, link added
$ws.Range("J9").Value = "snow_angels.jpg"
$ws.Range("L9").Value = "https://unsplash.com/photos/_hdO_l751fE"

# Row 10: snowman -- picture changed, link added
$ws.Range("J10").Value = "snowman.png"
$ws.Range("L10").Value = "https://unsplash.com/photos/5WIqleHzOok"

# Row 11: hot_chocolate_inside -- picture changed, link added
$ws.Range("J11").Value = "hot_chocolate_inside.jpg"
$ws.Range("L11").Value = "https://unsplash.com/s/photos/cozy-inside"

# Row 12: puzzle_inside -- picture changed, description + link added
$ws.Range("J12").Value = "puzzle_inside.jpg"
$ws.Range("K12").Value = "Make a puzzle inside"
$ws.Range("L12").Value = "https://unsplash.com/photos/AoX_1zm1NOM"

# Row 13: biking -- picture changed, description + link added
$ws.Range("J13").Value = "biking.jpg"
$ws.Range("K13").Value = "Go biking!"
$ws.Range("L13").Value = "https://unsplash.com/photos/JOnaeVoNkTQ"

# Row 14: sidewalk_drawing -- picture changed, description + link added
$ws.Range("J14").Value = "sidewalk_drawing.jpg"
$ws.Range("K14").Value = "Draw on the sidewalk!"
$ws.Range("L14").Value = "https://unsplash.com/photos/erxT2em063k"

# Row 15: water_gun_fight -- picture changed, description + link added
$ws.Range("J15").Value = "water_gun_fight.jpg"
$ws.Range("K15").Value = "Play with water guns"
$ws.Range("L15").Value = "https://unsplash.com/photos/JSLuw23jedY"

# Row 16: bake_cookies -- picture changed, description + link added
$ws.Range("J16").Value = "bake_cookies.jpg"
$ws.Range("K16").Value = "Bake cookies or cakes"
$ws.Range("L16").Value = "https://unsplash.com/photos/UyEmagArOLY"

# New row 29: read_inside activity
$ws.Range("A29").Value = "read_inside"
$ws.Range("B29").Value = -100
$ws.Range("C29").Value = 5
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 100
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 1000
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 100

# Update selection to match the committed state
[void]$ws.Range("J16").Select()
